$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 51 (pushes the existing rows 51-55 down to 52-56),
# then clear any formatting it inherited from the row above (e.g. hyperlink style)
# and fill in the new "Peak Design Leash" wishlist entry.
$ws.Rows("51").Insert()
$ws.Range("A51:D51").ClearFormats()

$ws.Range("C51").Value = "https://www.peakdesign.com/eu/products/leash?Color=Ibis"
$ws.Range("D51").Value = "50 EUR"
$ws.Range("B51").Value = "https://cdn.shopify.com/s/files/1/2986/1172/files/Leash_TopDown_Ibis.jpg?v=1753733022&width=700&height=700&crop=center"
$ws.Range("A51").Value = "Peak Design Leash"

# The autofilter / filter-database range needs to grow from A1:E55 to A1:E56 to
# include the newly inserted row, while keeping the existing "Y"-or-blank filter
# on column E (the 5th column of the range).
$ws.AutoFilterMode = $false
$critera = @("Y", "")
$ws.Range("A1:E56").AutoFilter(5, $critera, 7)

# Row insertion does not automatically grow the workbook-level hidden
# "_FilterDatabase" defined name that backs the autofilter - update it explicitly
# so it keeps pointing at the full A1:E56 range.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Wishlist Valentin!_FilterDatabase") {
        $n.RefersTo = "='Wishlist Valentin'!`$A`$1:`$E`$56"
    }
}

# Reflect the cursor position left behind in the saved workbook.
$ws.Range("A63").Select()
